# This script reproduces the "add new results to be analysed" commit:
# it extends the results table on Sheet1 from columns B:U to B:AO by
# appending a second copy of the angle values (row 3) and a fresh set of
# 0/1 flags (row 4) in the new columns V:AO, and updates the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: angle values (radians), columns V through AO ---
$ws.Range("V3").Value = 0.31415926535897898
$ws.Range("W3").Value = 0.62831853071795896
$ws.Range("X3").Value = 0.94247779607693805
$ws.Range("Y3").Value = 1.2566370614359199
$ws.Range("Z3").Value = 1.5707963267949001
$ws.Range("AA3").Value = 1.8849555921538801
$ws.Range("AB3").Value = 2.1991148575128601
$ws.Range("AC3").Value = 2.5132741228718301
$ws.Range("AD3").Value = 2.8274333882308098
$ws.Range("AE3").Value = 3.14159265358979
$ws.Range("AF3").Value = 3.4557519189487702
$ws.Range("AG3").Value = 3.76991118430775
$ws.Range("AH3").Value = 4.0840704496667302
$ws.Range("AI3").Value = 4.3982297150257104
$ws.Range("AJ3").Value = 4.7123889803846897
$ws.Range("AK3").Value = 5.0265482457436699
$ws.Range("AL3").Value = 5.3407075111026501
$ws.Range("AM3").Value = 5.6548667764616303
$ws.Range("AN3").Value = 5.9690260418206096
$ws.Range("AO3").Value = 6.2831853071795898

# --- Row 4: binary flags (0/1), columns V through AO ---
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 1
$ws.Range("Z4").Value = 1
$ws.Range("AA4").Value = 1
$ws.Range("AB4").Value = 1
$ws.Range("AC4").Value = 1
$ws.Range("AD4").Value = 1
$ws.Range("AE4").Value = 1
$ws.Range("AF4").Value = 1
$ws.Range("AG4").Value = 1
$ws.Range("AH4").Value = 1
$ws.Range("AI4").Value = 1
$ws.Range("AJ4").Value = 1
$ws.Range("AK4").Value = 1
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 1
$ws.Range("AO4").Value = 1

# --- View state: scroll the sheet right and move the selection ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$ws.Range("Z8").Select()

# --- Workbook window position (best effort) ---
$win.Left = 1155
$win.Top = 3270

